$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.819.90'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.46%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.544.44'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '582.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.39'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -3.95%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.583'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.45%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.105'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -3.39%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '5.56'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.44%  '
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.77%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '26.86'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.38%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.999.11'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '62.667.92'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.51%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000144'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.533.76'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.98'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -3.31%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '338.37'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -1.37%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.27'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.64%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.56'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.44%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '66.99'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.70%  '
$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.55'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.28%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.58'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.94%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.163'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -4.13%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.88'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.21%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.10'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.93%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.93'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.29%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '463.46'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0₃0787'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.77%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.65'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '175.72'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.395'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.61%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '18.67'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.62%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.46'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.72%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.09%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.69'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.63%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '40.02'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.46%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '156.62'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.65%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.66'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.40%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.04'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.625'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.95%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0530'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.50%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0954'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.26%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0234'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '17.90'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.94%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '11.36'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.65'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -5.61%  '
